# Auto-generated edit script: updates cached market-price/profit
# figures across multiple sheets to match the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ---------------- Sheet: ALC ----------------
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 3975
$ws.Range("I15").Value = 3975
$ws.Range("K15").Value = 11925
$ws.Range("M15").Value = -11756

# Row 28
$ws.Range("H28").Value = 283.86667
$ws.Range("I28").Value = 223.1
$ws.Range("K28").Value = 223.1
$ws.Range("M28").Value = 261.9

# Row 33
$ws.Range("H33").Value = 100.6
$ws.Range("I33").Value = 76.5
$ws.Range("K33").Value = 76.5
$ws.Range("M33").Value = 152.5

# Row 43
$ws.Range("H43").Value = 1275
$ws.Range("J43").Value = 1275
$ws.Range("L43").Value = 1275
$ws.Range("N43").Value = -1413

# Row 76
$ws.Range("H76").Value = 3353431
$ws.Range("I76").Value = 3910737
$ws.Range("J76").Value = 9595
$ws.Range("K76").Value = 3910737
$ws.Range("L76").Value = 9595
$ws.Range("M76").Value = -3910422
$ws.Range("N76").Value = -10225

# Row 79
$ws.Range("H79").Value = 3353431
$ws.Range("I79").Value = 3910737
$ws.Range("J79").Value = 9595
$ws.Range("K79").Value = 3910737
$ws.Range("L79").Value = 9595
$ws.Range("M79").Value = -3909645
$ws.Range("N79").Value = -11779

# Row 88
$ws.Range("H88").Value = 2818.0908
$ws.Range("J88").Value = 3062.5
$ws.Range("L88").Value = 3062.5
$ws.Range("N88").Value = -3874.5

# Row 91
$ws.Range("H91").Value = 2818.0908
$ws.Range("J91").Value = 3062.5
$ws.Range("L91").Value = 3062.5
$ws.Range("N91").Value = -5870.5

# Row 135
$ws.Range("H135").Value = 578.1053000000001
$ws.Range("I135").Value = 499.125
$ws.Range("J135").Value = 999.3333
$ws.Range("K135").Value = 4492.125
$ws.Range("L135").Value = 8993.9997
$ws.Range("M135").Value = -1957.125
$ws.Range("N135").Value = -14063.9997

# Row 138
$ws.Range("H138").Value = 3227.3157
$ws.Range("J138").Value = 2170.3438
$ws.Range("L138").Value = 6511.0314
$ws.Range("N138").Value = -16791.0314


# ---------------- Sheet: ARM ----------------
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2388.2727
$ws.Range("I32").Value = 1582.7632
$ws.Range("J32").Value = 7489.8335
$ws.Range("K32").Value = 1582.7632
$ws.Range("L32").Value = 7489.8335
$ws.Range("M32").Value = -1295.7632
$ws.Range("N32").Value = -8063.8335

# Row 61
$ws.Range("H61").Value = 2574.6365
$ws.Range("I61").Value = 1324.5
$ws.Range("K61").Value = 1324.5
$ws.Range("M61").Value = -1112.5

# Row 74
$ws.Range("H74").Value = 1596.2273
$ws.Range("I74").Value = 1467.3334
$ws.Range("K74").Value = 1467.3334
$ws.Range("M74").Value = -593.3334

# Row 77
$ws.Range("H77").Value = 1596.2273
$ws.Range("I77").Value = 1467.3334
$ws.Range("K77").Value = 7336.666999999999
$ws.Range("M77").Value = -2968.666999999999

# Row 110
$ws.Range("H110").Value = 3553.25
$ws.Range("I110").Value = 1400
$ws.Range("K110").Value = 1400
$ws.Range("M110").Value = 645

# Row 132
$ws.Range("H132").Value = 1415.804
$ws.Range("I132").Value = 1053.3684
$ws.Range("K132").Value = 3160.1052
$ws.Range("M132").Value = -630.1052

# Row 136
$ws.Range("H136").Value = 2574.6365
$ws.Range("I136").Value = 1324.5
$ws.Range("K136").Value = 3973.5
$ws.Range("M136").Value = -1423.5


# ---------------- Sheet: BSM ----------------
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1592.8
$ws.Range("I20").Value = 1795.5454
$ws.Range("K20").Value = 1795.5454
$ws.Range("M20").Value = -1548.5454

# Row 86
$ws.Range("H86").Value = 202891
$ws.Range("I86").Value = 2968.6667
$ws.Range("K86").Value = 2968.6667
$ws.Range("M86").Value = -1845.6667

# Row 89
$ws.Range("H89").Value = 202891
$ws.Range("I89").Value = 2968.6667
$ws.Range("K89").Value = 14843.3335
$ws.Range("M89").Value = -9227.333500000001

# Row 134
$ws.Range("H134").Value = 2014.4688
$ws.Range("I134").Value = 2019.4138
$ws.Range("K134").Value = 6058.2414
$ws.Range("M134").Value = -3523.2414


# ---------------- Sheet: CRP ----------------
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 750.5
$ws.Range("I16").Value = 750.5
$ws.Range("K16").Value = 750.5
$ws.Range("M16").Value = -463.5

# Row 31
$ws.Range("H31").Value = 2110.3635
$ws.Range("I31").Value = 1800
$ws.Range("J31").Value = 2482.8
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 2482.8
$ws.Range("M31").Value = -1505
$ws.Range("N31").Value = -3072.8

# Row 34
$ws.Range("H34").Value = 2110.3635
$ws.Range("I34").Value = 1800
$ws.Range("J34").Value = 2482.8
$ws.Range("K34").Value = 1800
$ws.Range("L34").Value = 2482.8
$ws.Range("M34").Value = -1598
$ws.Range("N34").Value = -2886.8

# Row 58
$ws.Range("H58").Value = 1550.1471
$ws.Range("I58").Value = 899.4583
$ws.Range("K58").Value = 899.4583
$ws.Range("M58").Value = -696.4583

# Row 62
$ws.Range("H62").Value = 2261.4443
$ws.Range("I62").Value = 2701.3333
$ws.Range("J62").Value = 2041.5
$ws.Range("K62").Value = 2701.3333
$ws.Range("L62").Value = 2041.5
$ws.Range("M62").Value = -2077.3333
$ws.Range("N62").Value = -3289.5

# Row 65
$ws.Range("H65").Value = 2261.4443
$ws.Range("I65").Value = 2701.3333
$ws.Range("J65").Value = 2041.5
$ws.Range("K65").Value = 13506.6665
$ws.Range("L65").Value = 10207.5
$ws.Range("M65").Value = -10386.6665
$ws.Range("N65").Value = -16447.5

# Row 113
$ws.Range("H113").Value = 750.5
$ws.Range("I113").Value = 750.5
$ws.Range("K113").Value = 750.5
$ws.Range("M113").Value = 1419.5

# Row 132
$ws.Range("H132").Value = 2444.7334
$ws.Range("I132").Value = 1560.5883
$ws.Range("K132").Value = 4681.7649
$ws.Range("M132").Value = -2151.7649

# Row 134
$ws.Range("H134").Value = 1686.7941
$ws.Range("I134").Value = 1448.4
$ws.Range("K134").Value = 4345.200000000001
$ws.Range("M134").Value = -1810.200000000001

# Row 136
$ws.Range("H136").Value = 1550.1471
$ws.Range("I136").Value = 899.4583
$ws.Range("K136").Value = 2698.3749
$ws.Range("M136").Value = -148.3748999999998


# ---------------- Sheet: CUL ----------------
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 1810
$ws.Range("I3").Value = 734
$ws.Range("K3").Value = 2202
$ws.Range("M3").Value = -2090

# Row 36
$ws.Range("H36").Value = 950.25
$ws.Range("J36").Value = 149
$ws.Range("L36").Value = 447
$ws.Range("N36").Value = -785

# Row 131
$ws.Range("H131").Value = 793.83
$ws.Range("J131").Value = 811.2447
$ws.Range("L131").Value = 2433.7341
$ws.Range("N131").Value = -12513.7341


# ---------------- Sheet: LTW ----------------
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4345.1875
$ws.Range("I40").Value = 1868.1111
$ws.Range("K40").Value = 1868.1111
$ws.Range("M40").Value = -1732.1111

# Row 61
$ws.Range("H61").Value = 3486.625
$ws.Range("I61").Value = 3379.6
$ws.Range("J61").Value = 3665
$ws.Range("K61").Value = 3379.6
$ws.Range("L61").Value = 3665
$ws.Range("M61").Value = -3177.6
$ws.Range("N61").Value = -4069

# Row 68
$ws.Range("H68").Value = 1950
$ws.Range("I68").Value = 1950
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1201
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 1950
$ws.Range("I71").Value = 1950
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9750
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6006
$ws.Range("N71").ClearContents()

# Row 110
$ws.Range("H110").Value = 36644
$ws.Range("J110").Value = 36644
$ws.Range("L110").Value = 36644
$ws.Range("N110").Value = -44824

# Row 113
$ws.Range("H113").Value = 3486.625
$ws.Range("I113").Value = 3379.6
$ws.Range("J113").Value = 3665
$ws.Range("K113").Value = 3379.6
$ws.Range("L113").Value = 3665
$ws.Range("M113").Value = -1209.6
$ws.Range("N113").Value = -8005

# Row 132
$ws.Range("H132").Value = 2447.0293
$ws.Range("I132").Value = 2073.3572
$ws.Range("K132").Value = 6220.071599999999
$ws.Range("M132").Value = -3690.071599999999


# ---------------- Sheet: WVR ----------------
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1331.3611
$ws.Range("I132").Value = 1046.9615
$ws.Range("J132").Value = 2070.8
$ws.Range("K132").Value = 3140.8845
$ws.Range("L132").Value = 6212.400000000001
$ws.Range("M132").Value = -610.8844999999997
$ws.Range("N132").Value = -11272.4
